$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 640, shifting existing rows 640-670 down to 641-671
$ws.Rows.Item(640).Insert()

# Populate the newly inserted row 640 with the new data record
$ws.Cells.Item(640, 1).Value = 3
$ws.Cells.Item(640, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(640, 3).Value = "Coquimbo"
$ws.Cells.Item(640, 4).Value = 45267
$ws.Cells.Item(640, 5).Value = 5
$ws.Cells.Item(640, 6).Value = 100114013
$ws.Cells.Item(640, 7).Value = "Zanahoria"
$ws.Cells.Item(640, 8).Value = "Sin especificar"
$ws.Cells.Item(640, 9).Value = "Primera"
$ws.Cells.Item(640, 10).Value = 340
$ws.Cells.Item(640, 11).Value = 5500
$ws.Cells.Item(640, 12).Value = 6000
$ws.Cells.Item(640, 13).Value = 5735
$ws.Cells.Item(640, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(640, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(640, 16).Value = 287
$ws.Cells.Item(640, 17).Value = 20
$ws.Cells.Item(640, 18).Value = "Hortaliza"
